$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.563.79"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "2.998.65"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.55"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.00"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "2.997.75"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("E12").Value = "  +4.50%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.29"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "3.489.13"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.02"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "61.528.81"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "2.994.18"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.16"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.35"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.00"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -4.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.58"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.95"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.52"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "0.0₃0834"
$ws.Range("E35").Value = "  +5.18%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.24"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.33"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("E41").Value = "  +9.45%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "395.34"
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.60"
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").Value = "2.716.44"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.63"
$ws.Range("E48").Value = "  +3.24%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +1.63%  "
